$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("A10").Value = 45795
$ws.Range("A10").NumberFormat = "YYYY-MM-DD"
$ws.Range("B10").Value = "Pengeluaran"
$ws.Range("C10").Value = 800000
$ws.Range("D10").Value = "Pupuk 50 liter"
